$wb = $excel.ActiveWorkbook

# --- Users sheet: selection change only ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B18:C22").Select()

# --- Member sheet: insert First_Name, Last_Name columns (B:C) ---
$wsMember = $wb.Worksheets.Item("Member")
$wsMember.Columns("B:C").Insert()
$wsMember.Range("B1").Value = "First_Name"
$wsMember.Range("C1").Value = "Last_Name"
$wsMember.Range("B2").Value = "string"
$wsMember.Range("C2").Value = "string"
$wsMember.Range("B3").Value = "yash"
$wsMember.Range("C3").Value = "kumar"
$wsMember.Range("B4").Value = "abhishek "
$wsMember.Range("C4").Value = "kumar"
$wsMember.Range("B5").Value = "vidhu"
$wsMember.Range("C5").Value = "prakash"
$wsMember.Range("B6").Value = "saurabh"
$wsMember.Range("C6").Value = "kant"
$wsMember.Range("B7").Value = "tannu"
$wsMember.Range("C7").Value = "baghel"
$wsMember.Range("B8").Value = "okasha"
$wsMember.Range("C8").Value = "anjum"
$wsMember.Range("B9").Value = "prince"
$wsMember.Range("C9").Value = "kumar"
$wsMember.Range("B10").Value = "anshul"
$wsMember.Range("C10").Value = "kumar"
$wsMember.Range("B11").Value = "gaurav"
$wsMember.Range("C11").Value = "kumar"
$wsMember.Range("B12").Value = "aaditya"
$wsMember.Range("C12").Value = "kumar"
$wsMember.Range("B1:C12").Select()

# --- Register sheet: selection change only ---
$wsRegister = $wb.Worksheets.Item("Register")
$wsRegister.Range("C1").Select()

# --- Email_IDs sheet: row height changes for rows 12-23 ---
$wsEmail = $wb.Worksheets.Item("Email_IDs")
$wsEmail.Rows("12:23").RowHeight = 28.8

# --- Personal_Instructor sheet: insert First_Name, Last_Name columns (C:D) ---
$wsPI = $wb.Worksheets.Item("Personal_Instructor")
$wsPI.Columns("C:D").Insert()
$wsPI.Range("C1").Value = "First_Name"
$wsPI.Range("D1").Value = "Last_Name"
$wsPI.Range("C2").Value = "string"
$wsPI.Range("D2").Value = "string"
$wsPI.Range("C3").Value = "aniket"
$wsPI.Range("D3").Value = "jadav"
$wsPI.Range("C4").Value = "ankit"
$wsPI.Range("D4").Value = "raj "
$wsPI.Range("C5").Value = "ahswin"
$wsPI.Range("D5").Value = "jha"
$wsPI.Range("C6").Value = "bipul"
$wsPI.Range("D6").Value = "bharti"
$wsPI.Range("C7").Value = "deepak"
$wsPI.Range("D7").Value = "rawal"
$wsPI.Range("C8").Value = "divyank"
$wsPI.Range("D8").Value = "saner"
$wsPI.Range("C9").Value = "gagan"
$wsPI.Range("D9").Value = "choudhary "
$wsPI.Range("C10").Value = "golu"
$wsPI.Range("D10").Value = "singh "
$wsPI.Range("C11").Value = "harsh"
$wsPI.Range("D11").Value = "baghel "
$wsPI.Range("C12").Value = "rishi"
$wsPI.Range("D12").Value = "kumar"
$wsPI.Range("C1:D12").Select()
